$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2-10) down by one row into rows (3-11),
# preserving each row's B:G metric values. Column A (quarter label)
# stays put since it is keyed off the row's existing period label.
for ($r = 10; $r -ge 2; $r--) {
    for ($c = 2; $c -le 7; $c++) {
        $val = $ws.Cells.Item($r, $c).Value2
        $ws.Cells.Item($r + 1, $c).Value2 = $val
    }
}

# New values for row 2 (freshly computed period).
$newRow2 = @(0.1724578193461484, 0.39058239716261, 0.3033305724894426, 0.550754548314803, 0.5414156770869448, 15)
for ($c = 2; $c -le 7; $c++) {
    $ws.Cells.Item(2, $c).Value2 = $newRow2[$c - 2]
}
